$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "f"
$ws.Range("A2").Value = "b"
$ws.Range("B2").Value = "bb"
$ws.Range("B3").Value = "bc"
$ws.Range("A3").Value = "c"
$ws.Range("A4").Value = "d"
$ws.Range("B4").Value = "bd"
$ws.Range("C1").Value = "hfghg"
$ws.Range("D1").Value = 2321

[void]$ws.Range("D1").Select()
